# Update the "Förändrad" (Changed) date column (C) for rows 2-28
# from serial date 45544 (2024-09-09) to 45545 (2024-09-10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45544) {
        $cell.Value = 45545
    }
}
